$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jobtable")

# Clear the old B1 value (the new layout no longer uses B1)
$ws.Range("B1").ClearContents()

# Row 1: opening bracket pieces + the big CONCATENATE formula
$ws.Range("D1").Value = "["
$ws.Range("A1").Value = "[`""
$ws.Range("F1").Formula = "=CONCATENATE(D1,D2,D3,D4,D5,D6,D7,D8,D9,D10,D11,D12,D13,D14,D15,D16,D17,D18,D19,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51,D52,D53,D54,D55,D56,D57,D58,D59,D60,D61,D62,D63,D64,D65,D66,D67,D68,D69,D70,D71,D72,D73,D74,D75,D76,D77,D78,D79,D80,D81,D82,D83,D84,D85,D86,D87,D88,D89,D90,D91,D92,D93,D94,D95,D96,D97,D98,D99,D100)"

# Rows 2..99: number in A, separator in B, CONCATENATE formula in D
for ($row = 2; $row -le 99; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 1).Value = $n
    $ws.Cells.Item($row, 2).Value = "`", `""
    $ws.Cells.Item($row, 4).Formula = "=CONCATENATE(A$row,B$row)"
}

# Row 100: last number/separator, but D100 is a literal closing value (not a formula)
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = "`", `""
$ws.Cells.Item(100, 4).Value = "99`"]"

$ws.Range("F1").Select()
